# Weekly update: insert a new week's worth of Coliflor price records
# (date 2022-07-27, serial 44769) at the top of the existing date-ordered
# block (rows 739-742), shifting the previously existing rows 739-816
# down to 743-820.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows at 739, pushing rows 739:816 down to 743:820.
$ws.Range("A739:A742").EntireRow.Insert()

# Fill the 4 newly inserted rows with the new week's data.
# Columns: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg,
#          F Categoria ID, G Categoria, H Variedad, I Calidad, J Volumen,
#          K Precio minimo, L Precio maximo, M Precio promedio ponderado,
#          N Unidad de comercializacion, O Origen, P Precio $/Kg,
#          Q Kg o Unidades, R Clasificacion

$newRows = @(
    @{ Row = 739; D = 44769; I = "Primera"; J = 6600; K = 600;  L = 700;  M = 648; O = "Región Metropolitana" },
    @{ Row = 740; D = 44769; I = "Primera"; J = 2200; K = 600;  L = 600;  M = 600; O = "Región de O'Higgins" },
    @{ Row = 741; D = 44769; I = "Segunda"; J = 3200; K = 500;  L = 500;  M = 500; O = "Región Metropolitana" },
    @{ Row = 742; D = 44769; I = "Segunda"; J = 1200; K = 400;  L = 400;  M = 400; O = "Región de O'Higgins" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 6
    $ws.Cells.Item($row, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = 100112008
    $ws.Cells.Item($row, 7).Value = "Coliflor"
    $ws.Cells.Item($row, 8).Value = "Sin especificar"
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = "`$/unidad"
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.M
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
